{"js": "// Adicionado Info Note Patrick\n//\n// Fill in the missing numeric specs for the \"Notebook Megaware\" bullet\n// list:\n//   \"GB de mem\u00f3ria RAM;\"             -> \"8GB de mem\u00f3ria RAM;\"\n//   \"GB de HD;\"                      -> \"500GB de HD;\"\n//   \"Processador Intel;\"             -> \"Processador Intel i5;\"\n//   \"Sistema Operacional Windows 7;\" -> \"Sistema Operacional Windows 7 x64Bit;\"\n//\n// The search text for the first item (\"GB de mem\") also occurs (as a\n// substring) inside the \"Notebook Lenovo\" / \"Macbook Air\" (\"4GB de\n// mem\u00f3ria RAM;\") and \"16GB de mem\u00f3ria\" items, and \"GB de HD;\" also\n// occurs inside the already-complete \"500GB de HD;\" item, so every\n// match is disambiguated by looking at the full text of its containing\n// paragraph before editing it. insertText(..., \"Replace\") on the\n// search-result range only rewrites the text of the run(s) actually\n// covered by the match, leaving neighboring runs (e.g. the accented\n// \"\u00f3\" / \"ria RAM;\" runs that follow \"GB de mem\") untouched, exactly\n// like the target edit.\n\nasync function replaceInParagraph(body, searchText, expectedParagraphText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  const paras = [];\n  for (let i = 0; i < results.items.length; i++) {\n    paras.push(results.items[i].paragraphs.getFirst());\n  }\n  for (let i = 0; i < paras.length; i++) {\n    paras[i].load(\"text\");\n  }\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    if (paras[i].text === expectedParagraphText) {\n      results.items[i].insertText(newText, \"Replace\");\n      await context.sync();\n      return true;\n    }\n  }\n  return false;\n}\n\nconst body = context.document.body;\n\nawait replaceInParagraph(body, \"GB de mem\", \"GB de mem\\u00F3ria RAM;\", \"8GB de mem\");\nawait replaceInParagraph(body, \"GB de HD;\", \"GB de HD;\", \"500GB de HD;\");\nawait replaceInParagraph(body, \"Processador Intel;\", \"Processador Intel;\", \"Processador Intel i5;\");\nawait replaceInParagraph(body, \"Sistema Operacional Windows 7;\", \"Sistema Operacional Windows 7;\", \"Sistema Operacional Windows 7 x64Bit;\");\n", "ps1": "# Adicionado Info Note Patrick\n# Fill in the missing numeric specs for the \"Notebook Megaware\" bullet list:\n#   \"GB de mem\u00f3ria RAM;\"            -> \"8GB de mem\u00f3ria RAM;\"\n#   \"GB de HD;\"                     -> \"500GB de HD;\"\n#   \"Processador Intel;\"            -> \"Processador Intel i5;\"\n#   \"Sistema Operacional Windows 7;\"-> \"Sistema Operacional Windows 7 x64Bit;\"\n#\n# Each paragraph below is located by its exact current text (paragraph-mark\n# included) so the edit lands only on the intended \"Notebook Megaware\" list\n# items and not on the similarly-worded \"Notebook Lenovo\" / \"Macbook Air\" /\n# \"16GB de mem\u00f3ria\" entries elsewhere in the document.\n\n$d = $word.ActiveDocument\n\n$memoria = [char]243  # the accented \"\u00f3\" in \"mem\u00f3ria\"\n\n$targets = @(\n    @{ Match = \"GB de mem\" + $memoria + \"ria RAM;\" + [char]13; Find = \"GB de mem\"; Replace = \"8GB de mem\" },\n    @{ Match = \"GB de HD;\" + [char]13; Find = \"GB de HD;\"; Replace = \"500GB de HD;\" },\n    @{ Match = \"Processador Intel;\" + [char]13; Find = \"Processador Intel;\"; Replace = \"Processador Intel i5;\" },\n    @{ Match = \"Sistema Operacional Windows 7;\" + [char]13; Find = \"Sistema Operacional Windows 7;\"; Replace = \"Sistema Operacional Windows 7 x64Bit;\" }\n)\n\nforeach ($target in $targets) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -eq $target.Match) {\n            $r = $p.Range\n            $r.Find.ClearFormatting()\n            $r.Find.Text = $target.Find\n            $r.Find.Replacement.ClearFormatting()\n            $r.Find.Replacement.Text = $target.Replace\n            $r.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n            break\n        }\n    }\n}\n"}
